# Weekly refresh of the "Feria Lagunitas de Puerto Montt - Caqui" price sheet.
# The date (Fecha, col D) for each daily record has shifted to a new reporting
# date; the volume (Volumen, col M) and, for row 6/8, the commercialisation
# unit details (Unidad de comercialización / Precio $/Kg / Kg por unidad,
# cols Q/S/T) moved together with the rows that now carry the "14 kilos
# granel" packaging record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44330

# Row 3
$ws.Range("D3").Value = 44306
$ws.Range("M3").Value = 80

# Row 4
$ws.Range("D4").Value = 44313
$ws.Range("M4").Value = 120

# Row 5
$ws.Range("D5").Value = 44323
$ws.Range("M5").Value = 80

# Row 6
$ws.Range("D6").Value = 44316
$ws.Range("M6").Value = 120
$ws.Range("Q6").Value = "$/caja 10 kilos empedrada"
$ws.Range("S6").Value = 11500
$ws.Range("T6").Value = 1

# Row 7
$ws.Range("D7").Value = 44302
$ws.Range("M7").Value = 80

# Row 8
$ws.Range("D8").Value = 44309
$ws.Range("Q8").Value = "$/caja 14 kilos granel"
$ws.Range("S8").Value = 821
$ws.Range("T8").Value = 14

# Row 9
$ws.Range("D9").Value = 44327
$ws.Range("M9").Value = 60

# Row 10
$ws.Range("D10").Value = 44322
$ws.Range("M10").Value = 60
